$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text (shared string) changes ---
$ws.Range("A3").Value = "1-Jul-2024 to 19-Dec-2024"
$ws.Range("B6").Value = "1-Jul-2024 to 19-Dec-2024"
$ws.Range("A160").Value = "4221 PATRIKA (AA GYA)"
$ws.Range("A299").Value = "5909 PATRIKA {M} (AA GYA)"
$ws.Range("A314").Value = "5924 PATRIKA {M} (400)"
$ws.Range("A369").Value = "6531 PATRIKA {M} (1500)"
$ws.Range("A448").Value = "7310 PATRIKA (D) (AA GYA)"

# --- Numeric value changes ---
$ws.Range("B10").Value = 314
$ws.Range("D10").Value = 314
$ws.Range("B11").Value = 228
$ws.Range("D11").Value = 209.82
$ws.Range("B20").Value = 129.5
$ws.Range("D20").Value = 259
$ws.Range("B23").Value = 119.5
$ws.Range("D23").Value = 239
$ws.Range("B24").Value = 60.5
$ws.Range("D24").Value = 127.05
$ws.Range("B50").Value = 129
$ws.Range("D50").Value = 141.9
$ws.Range("B54").Value = 30
$ws.Range("D54").Value = 33
$ws.Range("B72").Value = 489
$ws.Range("D72").Value = 757.95
$ws.Range("B74").Value = 700.5
$ws.Range("D74").Value = 1085.78
$ws.Range("B75").Value = 160.5
$ws.Range("B78").Value = 64.5
$ws.Range("B80").Value = 150
$ws.Range("B83").Value = 709
$ws.Range("D83").Value = 992.6
$ws.Range("B90").Value = 511.5
$ws.Range("D90").Value = 644.49
$ws.Range("B110").Value = 18.5
$ws.Range("D110").Value = 231.25
$ws.Range("B131").Value = 8
$ws.Range("D131").Value = 146.96
$ws.Range("B141").Value = 49
$ws.Range("D141").Value = 137.19999999999999
$ws.Range("B142").Value = 107
$ws.Range("D142").Value = 294.25
$ws.Range("B144").Value = 32
$ws.Range("D144").Value = 150.4
$ws.Range("B145").Value = 76
$ws.Range("D145").Value = 258.39999999999998
$ws.Range("B149").Value = 99.5
$ws.Range("D149").Value = 278.60000000000002
$ws.Range("B165").Value = 31
$ws.Range("D165").Value = 103.23
$ws.Range("B177").Value = 30.5
$ws.Range("D177").Value = 137.25
$ws.Range("B179").Value = 123.5
$ws.Range("D179").Value = 485.37
$ws.Range("B180").Value = 95.5
$ws.Range("D180").Value = 376.51
$ws.Range("B197").Value = 14.5
$ws.Range("D197").Value = 55.1
$ws.Range("B201").Value = 26.5
$ws.Range("D201").Value = 108.91
$ws.Range("B202").Value = 36
$ws.Range("D202").Value = 148.18
$ws.Range("B210").Value = 290.5
$ws.Range("D210").Value = 1307.25
$ws.Range("B212").Value = 25.5
$ws.Range("D212").Value = 133.88
$ws.Range("B225").Value = 8.5
$ws.Range("D225").Value = 51
$ws.Range("B240").Value = 45.5
$ws.Range("D240").Value = 159.25
$ws.Range("B252").Value = 33
$ws.Range("D252").Value = 165
$ws.Range("B264").Value = 27.5
$ws.Range("D264").Value = 171.88
$ws.Range("B267").Value = 36.5
$ws.Range("D267").Value = 219
$ws.Range("B269").Value = 72.5
$ws.Range("D269").Value = 430.65
$ws.Range("B288").Value = 25
$ws.Range("D288").Value = 179.44
$ws.Range("B289").Value = 15.5
$ws.Range("D289").Value = 110.52
$ws.Range("B313").Value = 9
$ws.Range("D313").Value = 99
$ws.Range("B318").Value = 91.5
$ws.Range("D318").Value = 503.25
$ws.Range("B327").Value = 8.5
$ws.Range("D327").Value = 58.65
$ws.Range("B332").Value = 14
$ws.Range("D332").Value = 68.599999999999994
$ws.Range("B337").Value = 64
$ws.Range("D337").Value = 416
$ws.Range("B351").Value = 0.5
$ws.Range("D351").Value = 3.92
$ws.Range("B368").Value = 20.5
$ws.Range("D368").Value = 160.72
$ws.Range("B369").Value = 8.5
$ws.Range("D369").Value = 212.5
$ws.Range("B371").Value = 34.5
$ws.Range("D371").Value = 377.09
$ws.Range("B383").Value = 22
$ws.Range("D383").Value = 240.46
$ws.Range("B390").Value = 28
$ws.Range("D390").Value = 56
$ws.Range("B403").Value = 80
$ws.Range("D403").Value = 216
$ws.Range("B423").Value = 25
$ws.Range("D423").Value = 85
$ws.Range("B448").Value = 6.5
$ws.Range("D448").Value = 35.75
$ws.Range("B458").Value = 48
$ws.Range("D458").Value = 342.24
$ws.Range("B468").Value = 24
$ws.Range("D468").Value = 102.72
$ws.Range("B480").Value = 86.5
$ws.Range("D480").Value = 121.1
$ws.Range("B487").Value = 209.5
$ws.Range("D487").Value = 314.25
$ws.Range("B488").Value = 350.5
$ws.Range("D488").Value = 525.75
$ws.Range("B490").Value = 137
$ws.Range("D490").Value = 53.67
$ws.Range("B491").Value = 96
$ws.Range("D491").Value = 218.88
$ws.Range("B497").Value = 12
$ws.Range("D497").Value = 33.6
$ws.Range("B512").Value = 58
$ws.Range("D512").Value = 193.14
$ws.Range("B516").Value = 53.5
$ws.Range("D516").Value = 192.6
$ws.Range("B525").Value = 31.5
$ws.Range("D525").Value = 119.7
$ws.Range("B536").Value = 32.5
$ws.Range("D536").Value = 162.5
$ws.Range("B546").Value = 23
$ws.Range("D546").Value = 138
$ws.Range("B547").Value = 3
$ws.Range("D547").Value = 18
$ws.Range("B578").Value = 251
$ws.Range("D578").Value = 482.04
$ws.Range("B579").Value = -1.5
$ws.Range("D579").Value = -2.94
$ws.Range("B580").Value = 249
$ws.Range("D580").Value = 480.48
$ws.Range("B581").Value = 115.5
$ws.Range("D581").Value = 327.77
$ws.Range("B590").Value = 29.5
$ws.Range("B591").Value = 49.5
$ws.Range("B594").Value = 628.5
$ws.Range("D594").Value = 502.8
$ws.Range("B596").Value = 479.5
$ws.Range("D596").Value = 383.6
$ws.Range("B602").Value = 467.5
$ws.Range("D602").Value = 318.83999999999997
$ws.Range("B605").Value = 36240.47
$ws.Range("D605").Value = 103376.68

# --- Row 97: clear to blank style (B/C/D) ---
$ws.Range("B97").ClearContents()
$ws.Range("C97").ClearContents()
$ws.Range("D97").ClearContents()
$ws.Range("B97").NumberFormat = """""0"
$ws.Range("C97").NumberFormat = """""0"
$ws.Range("D97").NumberFormat = """""0"

# --- Row 98: fill with data style (B/C/D) ---
$ws.Range("B98").Value = 8
$ws.Range("B98").NumberFormat = """""0.00"" pcs"""
$ws.Range("C98").Value = 1.55
$ws.Range("C98").NumberFormat = """""0.00"
$ws.Range("D98").Value = 12.4
$ws.Range("D98").NumberFormat = """""0.00"

# --- Row 239: clear to blank style (B/C/D) ---
$ws.Range("B239").ClearContents()
$ws.Range("C239").ClearContents()
$ws.Range("D239").ClearContents()
$ws.Range("B239").NumberFormat = """""0"
$ws.Range("C239").NumberFormat = """""0"
$ws.Range("D239").NumberFormat = """""0"
